# Updates the cryptos price/volume table (columns D and E) for rows 2-51
# to the refreshed values captured in this run, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.701.02"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.211.65"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'291.99"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'86.35"
$ws.Range("E6").Value = "  +7.07%  "
$ws.Range("D7").Value = "'0.515"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").Value = "'30.29"
$ws.Range("E10").Value = "  +4.86%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "'6.33"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "2.550.47"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "2.211.63"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "'0.727"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").Value = "39.646.34"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "'11.40"
$ws.Range("E20").Value = "  +11.64%  "
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").Value = "'5.79"
$ws.Range("D23").Value = "'65.68"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").Value = "'235.76"
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "'22.70"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "'9.26"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").Value = "'32.72"
$ws.Range("E31").Value = "  +4.00%  "
$ws.Range("D32").Value = "'152.17"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D35").Value = "'0.0717"
$ws.Range("E35").Value = "  +3.87%  "
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "'2.79"
$ws.Range("E38").Value = "  +7.03%  "
$ws.Range("D39").Value = "'15.92"
$ws.Range("E39").Value = "  +5.14%  "
$ws.Range("D40").Value = "'0.0986"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("D42").Value = "2.076.67"
$ws.Range("E42").Value = "  +9.58%  "
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").Value = "'0.0268"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("D46").Value = "'9.97"
$ws.Range("E46").Value = "  +11.05%  "
$ws.Range("D47").Value = "'17.68"
$ws.Range("E47").Value = "  +11.74%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "2.422.18"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "'70.79"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'88.99"
$ws.Range("E51").Value = "  +2.54%  "
